# Natmi following Dr Hou advice
# Rebuild the LR-pairs data block (rows 2-17, columns A:T) on the active sheet
# with the updated per-cluster-pair NATMI statistics (ECs/FAPs/M2/sCs x ECs/FAPs/M2/sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "C3", "C3ar1", "ECs", 3.0, 1.0, 1.558564, 4.675692, 0.005692101168584756, 0.005692101168584756, 1.0, 0.3333333333333333, 0.321326, 0.963978, 0.001759408141451291, 0.001759408141451291, 0.5008071358639999, 4.507264222776, 0.00001001472913797243, 0.00001001472913797243),
    @("ECs", "C3", "C3ar1", "FAPs", 3.0, 1.0, 1.558564, 4.675692, 0.005692101168584756, 0.005692101168584756, 1.0, 0.3333333333333333, 0.1103976666666667, 0.331193, 0.0006044781733521692, 0.0006044781733521692, 0.1720618289506667, 1.548556460556, 0.000003440750916921861, 0.000003440750916921861),
    @("ECs", "C3", "C3ar1", "M2", 3.0, 1.0, 1.558564, 4.675692, 0.005692101168584756, 0.005692101168584756, 3.0, 1.0, 178.9428813333334, 536.828644, 0.9797948571625668, 0.9797948571625668, 278.8939329024054, 2510.045396121648, 0.00557709145142838, 0.00557709145142838),
    @("ECs", "C3", "C3ar1", "sCs", 3.0, 1.0, 1.558564, 4.675692, 0.005692101168584756, 0.005692101168584756, 3.0, 1.0, 3.258402333333333, 9.775207, 0.01784125652262982, 0.01784125652262982, 5.078428574249332, 45.70585716824399, 0.0001015542371014816, 0.0001015542371014816),
    @("FAPs", "C3", "C3ar1", "ECs", 3.0, 1.0, 259.5505726666667, 778.6517180000001, 0.9479162344201305, 0.9479162344201304, 1.0, 0.3333333333333333, 0.321326, 0.963978, 0.001759408141451291, 0.001759408141451291, 83.40034731268935, 750.6031258142041, 0.001667771540252628, 0.001667771540252628),
    @("FAPs", "C3", "C3ar1", "FAPs", 3.0, 1.0, 259.5505726666667, 778.6517180000001, 0.9479162344201305, 0.9479162344201304, 1.0, 0.3333333333333333, 0.1103976666666667, 0.331193, 0.0006044781733521692, 0.0006044781733521692, 28.65377760439712, 257.883998439574, 0.0005729946738731471, 0.000572994673873147),
    @("FAPs", "C3", "C3ar1", "M2", 3.0, 1.0, 259.5505726666667, 778.6517180000001, 0.9479162344201305, 0.9479162344201304, 3.0, 1.0, 178.9428813333334, 536.828644, 0.9797948571625668, 0.9797948571625668, 46444.72732469006, 418002.5459222105, 0.9287634515057499, 0.9287634515057498),
    @("FAPs", "C3", "C3ar1", "sCs", 3.0, 1.0, 259.5505726666667, 778.6517180000001, 0.9479162344201305, 0.9479162344201304, 3.0, 1.0, 3.258402333333333, 9.775207, 0.01784125652262982, 0.01784125652262982, 845.7201915950697, 7611.481724355626, 0.01691201670025485, 0.01691201670025485),
    @("M2", "C3", "C3ar1", "ECs", 3.0, 1.0, 12.18925266666667, 36.567758, 0.04451691386950307, 0.04451691386950307, 1.0, 0.3333333333333333, 0.321326, 0.963978, 0.001759408141451291, 0.001759408141451291, 3.916723802369333, 35.250514221324, 0.00007832342069428959, 0.00007832342069428959),
    @("M2", "C3", "C3ar1", "FAPs", 3.0, 1.0, 12.18925266666667, 36.567758, 0.04451691386950307, 0.04451691386950307, 1.0, 0.3333333333333333, 0.1103976666666667, 0.331193, 0.0006044781733521692, 0.0006044781733521692, 1.345665052810445, 12.110985475294, 0.00002690950277911306, 0.00002690950277911306),
    @("M2", "C3", "C3ar1", "M2", 3.0, 1.0, 12.18925266666667, 36.567758, 0.04451691386950307, 0.04451691386950307, 3.0, 1.0, 178.9428813333334, 536.828644, 0.9797948571625668, 0.9797948571625668, 2181.17999347335, 19630.61994126015, 0.04361744326608805, 0.04361744326608805),
    @("M2", "C3", "C3ar1", "sCs", 3.0, 1.0, 12.18925266666667, 36.567758, 0.04451691386950307, 0.04451691386950307, 3.0, 1.0, 3.258402333333333, 9.775207, 0.01784125652262982, 0.01784125652262982, 39.71748933065622, 357.457403975906, 0.0007942376799416215, 0.0007942376799416216),
    @("sCs", "C3", "C3ar1", "ECs", 3.0, 1.0, 0.5133286666666667, 1.539986, 0.001874750541781658, 0.001874750541781658, 1.0, 0.3333333333333333, 0.321326, 0.963978, 0.001759408141451291, 0.001759408141451291, 0.1649458471453333, 1.484512624308, 0.000003298451366400867, 0.000003298451366400867),
    @("sCs", "C3", "C3ar1", "FAPs", 3.0, 1.0, 0.5133286666666667, 1.539986, 0.001874750541781658, 0.001874750541781658, 1.0, 0.3333333333333333, 0.1103976666666667, 0.331193, 0.0006044781733521692, 0.0006044781733521692, 0.05667028703311111, 0.5100325832979999, 0.000001133245782987166, 0.000001133245782987166),
    @("sCs", "C3", "C3ar1", "M2", 3.0, 1.0, 0.5133286666666667, 1.539986, 0.001874750541781658, 0.001874750541781658, 3.0, 1.0, 178.9428813333334, 536.828644, 0.9797948571625668, 0.9797948571625668, 91.85651068433157, 826.708596158984, 0.001836870939300404, 0.001836870939300404),
    @("sCs", "C3", "C3ar1", "sCs", 3.0, 1.0, 0.5133286666666667, 1.539986, 0.001874750541781658, 0.001874750541781658, 3.0, 1.0, 3.258402333333333, 9.775207, 0.01784125652262982, 0.01784125652262982, 1.672631325233555, 15.053681927102, 0.00003344790533186579, 0.0000334479053318658)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $rowIndex = $i + 2
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}

Write-Host "Updated $($data.Length) data rows (A2:T17)"
